$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.749.50"
$ws.Range("E2").Value = "  +7.09%  "
$ws.Range("D3").Value = "3.469.62"
$ws.Range("E3").Value = "  +5.22%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "413.92"
$ws.Range("E5").Value = "  +3.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.75"
$ws.Range("E6").Value = "  +18.37%  "
$ws.Range("D7").Value = "3.461.27"
$ws.Range("E7").Value = "  +5.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.592"
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.693"
$ws.Range("E10").Value = "  +9.28%  "
$ws.Range("E11").Value = "  +29.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.46"
$ws.Range("E12").Value = "  +6.75%  "
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").Value = "4.019.49"
$ws.Range("E14").Value = "  +5.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.72"
$ws.Range("E15").Value = "  +4.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.10"
$ws.Range("E16").Value = "  +4.47%  "
$ws.Range("D17").Value = "3.461.51"
$ws.Range("E17").Value = "  +4.95%  "
$ws.Range("D18").Value = "62.617.48"
$ws.Range("E18").Value = "  +7.34%  "
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.81"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000138"
$ws.Range("E21").Value = "  +25.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.37"
$ws.Range("E22").Value = "  +1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "82.30"
$ws.Range("E23").Value = "  +10.35%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "316.41"
$ws.Range("E24").Value = "  +4.76%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.15"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.77"
$ws.Range("E27").Value = "  +8.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.18"
$ws.Range("E28").Value = "  +3.95%  "
$ws.Range("E29").Value = "  +4.14%  "
$ws.Range("B30").Value = "LEO"
$ws.Range("C30").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.39"
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.178"
$ws.Range("E31").Value = "  +4.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  +4.76%  "
$ws.Range("E33").Value = "  +24.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.77"
$ws.Range("E34").Value = "  +3.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "43.01"
$ws.Range("E35").Value = "  +5.44%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  -6.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.44"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.56"
$ws.Range("E39").Value = "  +2.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.996"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("E41").Value = "  -7.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.02"
$ws.Range("E42").Value = "  +7.35%  "
$ws.Range("E43").Value = "  +3.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "136.92"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.25"
$ws.Range("E45").Value = "  +2.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.287"
$ws.Range("E46").Value = "  +2.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.97"
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.15"
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("D50").Value = "2.227.19"
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("D51").Value = "3.821.82"
$ws.Range("E51").Value = "  +5.30%  "
